$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 323, shifting the existing row 323 (and
# everything below it) down by one. This matches the diff, where a brand
# new record appears at row 323 and all subsequent rows (old 323..414)
# move down to become rows 324..415.
$ws.Rows("323:323").Insert()

# Populate the newly inserted row 323 with the new record's data.
$ws.Range("A323").Value = 3
$ws.Range("B323").Value = "Femacal de La Calera"
$ws.Range("C323").Value = "Coquimbo"
$ws.Range("D323").Value = 44722
$ws.Range("E323").Value = 5
$ws.Range("F323").Value = "Fruta"
$ws.Range("G323").Value = 100108
$ws.Range("H323").Value = "Tropicales y subtropicales"
$ws.Range("I323").Value = 100108002
$ws.Range("J323").Value = "Mango"
$ws.Range("K323").Value = "Sin especificar"
$ws.Range("L323").Value = "Primera"
$ws.Range("M323").Value = 228
$ws.Range("N323").Value = 10000
$ws.Range("O323").Value = 10000
$ws.Range("P323").Value = 10000
$ws.Range("Q323").Value = "$/bandeja 4 kilos"
$ws.Range("R323").Value = "Brasil"
$ws.Range("S323").Value = 2500
$ws.Range("T323").Value = 4
